# "Generate Report for Handoff"
#
# The e8cedaa1-052c-49d2-afee-ad721cdcfa71.md and
# fb6f0c78-fe15-4b7b-8736-4299269edf03.md source files have moved on from
# "Handed back: in sync with en-US" to being ready for another handoff
# round, because the handback that was received is for a stale commit of
# the source file. Update the Overview sheet and both per-locale sheets
# (zh-cn, de-de) to reflect the new status/timestamps and record the
# "stale handback" error detail, then widen the Error Detail column so the
# new message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyStatus = "Ready for handoff"

$errDetailE8 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25070db70708d91cf6d8fb565cc9eaaff3bbe405/e2e/e8cedaa1-052c-49d2-afee-ad721cdcfa71.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3644f8931dc3662be1e67fa3c1abd947d88f3605/e2e/e8cedaa1-052c-49d2-afee-ad721cdcfa71.md."
$errDetailFb = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25070db70708d91cf6d8fb565cc9eaaff3bbe405/e2e/fb6f0c78-fe15-4b7b-8736-4299269edf03.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3644f8931dc3662be1e67fa3c1abd947d88f3605/e2e/fb6f0c78-fe15-4b7b-8736-4299269edf03.md."

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (e8cedaa1) and 5 (fb6f0c78) -> zh-cn (E) / de-de
# (F) status columns become "Ready for handoff" and the "Latest HO Xliff
# Generate Date" (G) column picks up the new generation timestamp.
# ---------------------------------------------------------------------
$wsOverview.Range("E4").Value = $readyStatus
$wsOverview.Range("F4").Value = $readyStatus
$wsOverview.Range("G4").Value = "2016-08-18 06:25:50"

$wsOverview.Range("E5").Value = $readyStatus
$wsOverview.Range("F5").Value = $readyStatus
$wsOverview.Range("G5").Value = "2016-08-18 06:25:50"

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4/5 Status -> Ready for handoff, Latest Handoff
# Datetime refreshed, Error Detail populated with the stale-handback note.
# ---------------------------------------------------------------------
$wsZhCn.Range("C4").Value = $readyStatus
$wsZhCn.Range("H4").Value = "2016-08-18 06:25:45"
$wsZhCn.Range("P4").Value = $errDetailE8

$wsZhCn.Range("C5").Value = $readyStatus
$wsZhCn.Range("H5").Value = "2016-08-18 06:25:45"
$wsZhCn.Range("P5").Value = $errDetailFb

# ---------------------------------------------------------------------
# de-de sheet: same treatment, different refreshed timestamp.
# ---------------------------------------------------------------------
$wsDeDe.Range("C4").Value = $readyStatus
$wsDeDe.Range("H4").Value = "2016-08-18 06:25:50"
$wsDeDe.Range("P4").Value = $errDetailE8

$wsDeDe.Range("C5").Value = $readyStatus
$wsDeDe.Range("H5").Value = "2016-08-18 06:25:50"
$wsDeDe.Range("P5").Value = $errDetailFb

# ---------------------------------------------------------------------
# Widen the Error Detail column (P, the 16th column) on both locale
# sheets so the long message is readable. ColumnWidth and the OOXML
# <col width> differ by the standard 5/6 character padding, so feed the
# setter 40 - 5/6 to land on a stored width of exactly 40.
# ---------------------------------------------------------------------
$targetColWidth = 40 - (5/6)
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
